$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = 1.33
$ws.Range("Q2").Value = 1.93
$ws.Range("S2").Value = 3.4
$ws.Range("G3").Value = 1.27
$ws.Range("H3").Value = 12.5
$ws.Range("I3").Value = 15
$ws.Range("K3").Value = 9.199999999999999
$ws.Range("N3").Value = 7.6
$ws.Range("Q3").Value = 1.26
$ws.Range("R3").Value = 2.06
$ws.Range("S3").Value = 1.65
$ws.Range("U3").Value = 2.02
$ws.Range("V3").Value = 1.07
$ws.Range("W3").Value = 4.7
$ws.Range("X3").Value = 65
$ws.Range("Y3").Value = 80
$ws.Range("Z3").Value = 180
$ws.Range("AA3").Value = 550
$ws.Range("AB3").Value = 19.5
$ws.Range("AC3").Value = 25
$ws.Range("AD3").Value = 60
$ws.Range("AE3").Value = 210
$ws.Range("AF3").Value = 14.5
$ws.Range("AG3").Value = 16
$ws.Range("AH3").Value = 36
$ws.Range("AI3").Value = 130
$ws.Range("AJ3").Value = 13.5
$ws.Range("AK3").Value = 16
$ws.Range("AL3").Value = 38
$ws.Range("AM3").Value = 130
$ws.Range("AN3").Value = 3.3
$ws.Range("AO3").Value = 180
$ws.Range("H4").Value = 2.22
$ws.Range("I4").Value = 2.24
$ws.Range("O4").Value = 1.46
$ws.Range("T4").Value = 2.06
$ws.Range("AE4").Value = 28
$ws.Range("AH4").Value = 23
$ws.Range("AJ4").Value = 110
$ws.Range("I5").Value = 2.96
$ws.Range("O5").Value = 1.43
$ws.Range("R5").Value = 1.26
$ws.Range("T5").Value = 1.96
$ws.Range("U5").Value = 1.98
$ws.Range("AE5").Value = 75
$ws.Range("G6").Value = 3.2
$ws.Range("I6").Value = 2.64
$ws.Range("S6").Value = 4
$ws.Range("Z6").Value = 16
$ws.Range("AA6").Value = 40
$ws.Range("AC6").Value = 7
$ws.Range("AH6").Value = 19.5
$ws.Range("AI6").Value = 46
$ws.Range("AJ6").Value = 55
$ws.Range("AL6").Value = 55
$ws.Range("AM6").Value = 130
